# Update Name of Algo
# Applies updated imputed values to result_data_KNN.xlsx (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.271000000000001
$ws.Range("B3").Value = 5.616999999999999
$ws.Range("B14").Value = 5.756
$ws.Range("B21").Value = 9.379000000000001
$ws.Range("B23").Value = 7.598000000000001
$ws.Range("B25").Value = 6.334000000000001
$ws.Range("D25").Value = -7.889
$ws.Range("B26").Value = 6.547
$ws.Range("D27").Value = -8.652999999999999
$ws.Range("B29").Value = 5.689
$ws.Range("D31").Value = -8.161000000000001
$ws.Range("D39").Value = -7.782999999999999
$ws.Range("D48").Value = -7.415000000000001
$ws.Range("D51").Value = -8.34
$ws.Range("D52").Value = -8.1
$ws.Range("B53").Value = 6.005000000000001
$ws.Range("D55").Value = -7.904000000000001
$ws.Range("D56").Value = -8.434999999999999
$ws.Range("B57").Value = 4.981999999999999
$ws.Range("D57").Value = -8.059999999999999
$ws.Range("B59").Value = 5.145
$ws.Range("B69").Value = 5.095
$ws.Range("D73").Value = -8.278000000000002
$ws.Range("B79").Value = 6.031
$ws.Range("B83").Value = 5.915
$ws.Range("D89").Value = -6.351
$ws.Range("D90").Value = -7.472
$ws.Range("B91").Value = 5.469
$ws.Range("D92").Value = -6.35
$ws.Range("B93").Value = 5.459

$wb.Save()
